$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows whose Target cluster (D) = "ECs" (old rows 8, 5, 2), bottom-up so indices stay valid
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# Update remaining rows (now rows 2-7) with recomputed TPM-based values
# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rspo3"
$ws.Range("C2").Value = "Lgr5"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.003420333333333333
$ws.Range("H2").Value = 0.010261
$ws.Range("I2").Value = 0.003549653112303053
$ws.Range("J2").Value = 0.003549653112303053
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.3961209999999999
$ws.Range("N2").Value = 1.188363
$ws.Range("O2").Value = 0.5646784620538419
$ws.Range("P2").Value = 0.5646784620538419
$ws.Range("Q2").Value = 0.001354865860333333
$ws.Range("R2").Value = 0.012193792743
$ws.Range("S2").Value = 0.002004412660279922
$ws.Range("T2").Value = 0.002004412660279922

# Row 3: ECs -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rspo3"
$ws.Range("C3").Value = "Lgr5"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.003420333333333333
$ws.Range("H3").Value = 0.010261
$ws.Range("I3").Value = 0.003549653112303053
$ws.Range("J3").Value = 0.003549653112303053
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 0.3053773333333333
$ws.Range("N3").Value = 0.916132
$ws.Range("O3").Value = 0.4353215379461581
$ws.Range("P3").Value = 0.4353215379461581
$ws.Range("Q3").Value = 0.001044492272444444
$ws.Range("R3").Value = 0.009400430452
$ws.Range("S3").Value = 0.001545240452023132
$ws.Range("T3").Value = 0.001545240452023132

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Lgr5"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3764723333333334
$ws.Range("H4").Value = 1.129417
$ws.Range("I4").Value = 0.3907064193682856
$ws.Range("J4").Value = 0.3907064193682855
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 0.3961209999999999
$ws.Range("N4").Value = 1.188363
$ws.Range("O4").Value = 0.5646784620538419
$ws.Range("P4").Value = 0.5646784620538419
$ws.Range("Q4").Value = 0.1491285971523333
$ws.Range("R4").Value = 1.342157374371
$ws.Range("S4").Value = 0.2206235000034469
$ws.Range("T4").Value = 0.2206235000034469

# Row 5: FAPs -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo3"
$ws.Range("C5").Value = "Lgr5"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2.0
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3764723333333334
$ws.Range("H5").Value = 1.129417
$ws.Range("I5").Value = 0.3907064193682856
$ws.Range("J5").Value = 0.3907064193682855
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 0.3053773333333333
$ws.Range("N5").Value = 0.916132
$ws.Range("O5").Value = 0.4353215379461581
$ws.Range("P5").Value = 0.4353215379461581
$ws.Range("Q5").Value = 0.1149661172271111
$ws.Range("R5").Value = 1.034695055044
$ws.Range("S5").Value = 0.1700829193648387
$ws.Range("T5").Value = 0.1700829193648387

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Rspo3"
$ws.Range("C6").Value = "Lgr5"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 0.5836756666666667
$ws.Range("H6").Value = 1.751027
$ws.Range("I6").Value = 0.6057439275194114
$ws.Range("J6").Value = 0.6057439275194113
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.3961209999999999
$ws.Range("N6").Value = 1.188363
$ws.Range("O6").Value = 0.5646784620538419
$ws.Range("P6").Value = 0.5646784620538419
$ws.Range("Q6").Value = 0.2312061887556666
$ws.Range("R6").Value = 2.080855698801
$ws.Range("S6").Value = 0.3420505493901151
$ws.Range("T6").Value = 0.342050549390115

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Rspo3"
$ws.Range("C7").Value = "Lgr5"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 0.5836756666666667
$ws.Range("H7").Value = 1.751027
$ws.Range("I7").Value = 0.6057439275194114
$ws.Range("J7").Value = 0.6057439275194113
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 0.3053773333333333
$ws.Range("N7").Value = 0.916132
$ws.Range("O7").Value = 0.4353215379461581
$ws.Range("P7").Value = 0.4353215379461581
$ws.Range("Q7").Value = 0.1782413186182222
$ws.Range("R7").Value = 1.604171867564
$ws.Range("S7").Value = 0.2636933781292963
$ws.Range("T7").Value = 0.2636933781292962
